$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the two table headers before we touch columns, since they live inside
# the range we are about to delete/reshuffle.
$i2plsValue = $ws.Range("B2").Value2
$kbtsValue = $ws.Range("G2").Value2

# Remove the extra spacer column (old F) and the duplicate KBTS "File" label column (old H),
# shifting the KBTS Profit/Time values left into F/G.
$ws.Columns("F:H").Delete()

# Move the "I2PLS" header from B2 (old merge start) to sit above the Profit/Time columns (D2:E2)
$ws.Range("B2:E2").UnMerge()
$ws.Range("B2").Value = $null
$ws.Range("D2").Value = $i2plsValue
$ws.Range("D2:E2").Merge()

# Restore the "KBTS" header above its Profit/Time columns (F2:G2)
$ws.Range("F2").Value = $kbtsValue

# Match the saved selection/cursor position recorded in the workbook
[void]$ws.Range("K8").Select()
